$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Deem753"
$ws.Range("B2").Value = 23070701
$ws.Range("C2").Value = "veer33"
$ws.Range("D2").Value = "aP$32r%K"
